# Daily scrape update - 2026-02-08 04:39:59 UTC
# Refreshes the Global Talent opportunity listing: new rows of scraped
# data replace the previous day's rows 2-9, three column widths are
# adjusted to fit the new content, and the stale "Yes/highlighted"
# premium flag on the old rows 3/4/9 is cleared (those opportunities are
# no longer premium in today's scrape).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's COM layer pads any ColumnWidth assignment by this constant
# (~5/6 of a character) before it lands in the saved <col width="..">
# attribute, so back it out to land on the exact target widths.
$pad = 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 53 - $pad
$ws.Columns.Item(4).ColumnWidth = 56 - $pad
$ws.Columns.Item(6).ColumnWidth = 17 - $pad
$ws.Columns.Item(8).ColumnWidth = 56 - $pad

function Set-Row($r, $id, $link, $title, $country, $premium, $applicants, $duration, $org) {
    # The opportunity id looks numeric, but the scrape always stores it as
    # plain text with no special cell style. Force text via a temporary
    # "@" number format so Excel doesn't auto-coerce it to a Number, then
    # clear formats again so the cell ends up with no style index at all
    # (matching the un-styled inline-string cells the scraper writes).
    $idCell = $ws.Range("A$r")
    $idCell.NumberFormat = "@"
    $idCell.Value = $id
    $idCell.ClearFormats()

    $ws.Range("B$r").Value = $link
    $ws.Range("C$r").Value = $title
    $ws.Range("D$r").Value = $country
    $ws.Range("E$r").Value = $premium
    $ws.Range("F$r").Value = $applicants
    $ws.Range("G$r").Value = $duration
    $ws.Range("H$r").Value = $org
}

Set-Row 2 "1331750" "https://aiesec.org/opportunity/global-talent/1331750" "Neuro-Marketing & Communications Intern" "Amman, Jordan" "No" "1 applicant" "9 - 12 Weeks" "Amoux Group"
Set-Row 3 "1331747" "https://aiesec.org/opportunity/global-talent/1331747" "Sales Representative" "Istanbul, İstanbul, Türkiye" "No" "5 applicants" "9 - 12 Weeks" "Apilex"
Set-Row 4 "1331473" "https://aiesec.org/opportunity/global-talent/1331473" "Business Adminstration" "Bengaluru, Karnataka, India" "No" "1 applicant" "9 - 12 Weeks" "The Perk Central Cafe"
Set-Row 5 "1330859" "https://aiesec.org/opportunity/global-talent/1330859" "Export Sales Specialist" "Başakşehir, Başak, 34490 Başakşehir/İstanbul, Türkiye" "No" "62 applicants" "9 - 12 Weeks" "Esen Isıtma Soğutma Elektrik İnşaat Sanayi ve Ticaret"
Set-Row 6 "1327286" "https://aiesec.org/opportunity/global-talent/1327286" "Taxes & Internal Control" "Panamá, Provincia de Panamá, Panamá" "No" "101 applicants" "6 - 18 Months" "NESTLÉ"
Set-Row 7 "1325297" "https://aiesec.org/opportunity/global-talent/1325297" "International Sales Representetive Spanish Speaker" "Maslak, Sarıyer/İstanbul, Türkiye" "No" "30 applicants" "6 - 18 Months" "Esvita Clinic"
Set-Row 8 "1321053" "https://aiesec.org/opportunity/global-talent/1321053" "International Sales Representetive German Speaker" "Maslak, Sarıyer/İstanbul, Türkiye" "No" "20 applicants" "6 - 18 Months" "Esvita Clinic"
Set-Row 9 "1321052" "https://aiesec.org/opportunity/global-talent/1321052" "International Sales Representetive" "Maslak, Sarıyer/İstanbul, Türkiye" "No" "165 applicants" "6 - 18 Months" "Esvita Clinic"

# Rows 3, 4 and 9 previously carried the yellow "premium" highlight
# (style s="3"); today's scrape no longer flags them, so drop the fill
# back to the sheet's default (unstyled) look.
$ws.Range("E3").ClearFormats()
$ws.Range("E4").ClearFormats()
$ws.Range("E9").ClearFormats()
